$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace target radius values: 16 -> 12, 11.314 -> 8.485 (and their negatives)
for ($r = 2; $r -le 41; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $aVal = $aCell.Value2
    $bVal = $bCell.Value2

    if ($aVal -eq 16) { $aCell.Value2 = 12 }
    elseif ($aVal -eq -16) { $aCell.Value2 = -12 }
    elseif ($aVal -eq 11.314) { $aCell.Value2 = 8.485 }
    elseif ($aVal -eq -11.314) { $aCell.Value2 = -8.485 }

    if ($bVal -eq 16) { $bCell.Value2 = 12 }
    elseif ($bVal -eq -16) { $bCell.Value2 = -12 }
    elseif ($bVal -eq 11.314) { $bCell.Value2 = 8.485 }
    elseif ($bVal -eq -11.314) { $bCell.Value2 = -8.485 }
}

# Update view: zoom and selection
$ws.Range("E18").Select()
$excel.ActiveWindow.Zoom = 115
